# Update SABR calibration outputs (alpha, rho, nu) on their respective sheets.
$wb = $excel.ActiveWorkbook

$alpha = $wb.Worksheets.Item("sabr_alpha")
$rho   = $wb.Worksheets.Item("sabr_rho")
$nu    = $wb.Worksheets.Item("sabr_nu")

$alphaVals = @(
    @(0.1391177694231066, 0.1846311437361156, 0.1968245453932313, 0.1781495538911682, 0.1727589651188632),
    @(0.1666405030017211, 0.1986920640478133, 0.1999451847614024, 0.1883856878566519, 0.17063106637228),
    @(0.1766681763817387, 0.1926930765646341, 0.1965753688622109, 0.1926050870219206, 0.1773786194049607)
)

$rhoVals = @(
    @(-0.6329876681687155, -0.525337597569106, -0.4831297395439832, -0.4165000751703339, -0.3080300238788374),
    @(-0.5843784817874683, -0.5451623032059372, -0.5192037559179212, -0.498743074072418, -0.3901202644965844),
    @(-0.5525023448268775, -0.5453249721517903, -0.5315120167547526, -0.5404873546671907, -0.5084018917440672)
)

$nuVals = @(
    @(2.049198633956934, 1.679633689423351, 1.441003835359067, 1.066331560365873, 0.7699940557453125),
    @(1.332680149730216, 1.064564699296813, 0.953271385366513, 0.6880500324821187, 0.5533080301251064),
    @(1.040785710879257, 0.9502018101951232, 0.8796305578792002, 0.7296450726108302, 0.6199948348583607)
)

for ($i = 0; $i -lt 3; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt 5; $j++) {
        $col = $j + 2
        $alpha.Cells.Item($row, $col).Value = $alphaVals[$i][$j]
        $rho.Cells.Item($row, $col).Value = $rhoVals[$i][$j]
        $nu.Cells.Item($row, $col).Value = $nuVals[$i][$j]
    }
}
